$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1066.5555
$ws.Range("I9").Value = 1016.8333
$ws.Range("J9").Value = 1166
$ws.Range("K9").Value = 1016.8333
$ws.Range("L9").Value = 1166
$ws.Range("M9").Value = -847.8333
$ws.Range("N9").Value = -1504

$ws.Range("H76").Value = 4331.6
$ws.Range("I76").Value = 3919.5
$ws.Range("J76").Value = 4949.75
$ws.Range("K76").Value = 3919.5
$ws.Range("L76").Value = 4949.75
$ws.Range("M76").Value = -3604.5
$ws.Range("N76").Value = -5579.75

$ws.Range("H79").Value = 4331.6
$ws.Range("I79").Value = 3919.5
$ws.Range("J79").Value = 4949.75
$ws.Range("K79").Value = 3919.5
$ws.Range("L79").Value = 4949.75
$ws.Range("M79").Value = -2827.5
$ws.Range("N79").Value = -7133.75

$ws.Range("H132").Value = 2394958.8
$ws.Range("I132").Value = 2660824.8
$ws.Range("K132").Value = 7982474.399999999
$ws.Range("M132").Value = -7979944.399999999

$ws.Range("H137").Value = 30499.908
$ws.Range("I137").Value = 32799.9
$ws.Range("K137").Value = 98399.70000000001
$ws.Range("M137").Value = -95849.70000000001

$ws.Range("H138").Value = 1955.091
$ws.Range("I138").Value = 1290.1063
$ws.Range("J138").Value = 3600.0527
$ws.Range("K138").Value = 3870.3189
$ws.Range("L138").Value = 10800.1581
$ws.Range("M138").Value = 1269.6811
$ws.Range("N138").Value = -21080.1581

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17366.3
$ws.Range("J32").Value = 5774.75
$ws.Range("L32").Value = 5774.75
$ws.Range("N32").Value = -6348.75

$ws.Range("H45").Value = 3173
$ws.Range("I45").Value = 1263.3334
$ws.Range("J45").Value = 3889.125
$ws.Range("K45").Value = 1263.3334
$ws.Range("L45").Value = 3889.125
$ws.Range("M45").Value = -886.3334
$ws.Range("N45").Value = -4643.125

$ws.Range("H61").Value = 10429.077
$ws.Range("I61").Value = 2612.5
$ws.Range("J61").Value = 17129
$ws.Range("K61").Value = 2612.5
$ws.Range("L61").Value = 17129
$ws.Range("M61").Value = -2400.5
$ws.Range("N61").Value = -17553

$ws.Range("H74").Value = 120547.22
$ws.Range("I74").Value = 128572.53
$ws.Range("K74").Value = 128572.53
$ws.Range("M74").Value = -127698.53

$ws.Range("H77").Value = 120547.22
$ws.Range("I77").Value = 128572.53
$ws.Range("K77").Value = 642862.65
$ws.Range("M77").Value = -638494.65

$ws.Range("H132").Value = 1953.1207
$ws.Range("I132").Value = 1781.262
$ws.Range("J132").Value = 2404.25
$ws.Range("K132").Value = 5343.786
$ws.Range("L132").Value = 7212.75
$ws.Range("M132").Value = -2813.786
$ws.Range("N132").Value = -12272.75

$ws.Range("H136").Value = 10429.077
$ws.Range("I136").Value = 2612.5
$ws.Range("J136").Value = 17129
$ws.Range("K136").Value = 7837.5
$ws.Range("L136").Value = 51387
$ws.Range("M136").Value = -5287.5
$ws.Range("N136").Value = -56487

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1266.875
$ws.Range("I64").Value = 1452.5
$ws.Range("J64").Value = 1081.25
$ws.Range("K64").Value = 1452.5
$ws.Range("L64").Value = 1081.25
$ws.Range("M64").Value = -1227.5
$ws.Range("N64").Value = -1531.25

$ws.Range("H67").Value = 1266.875
$ws.Range("I67").Value = 1452.5
$ws.Range("J67").Value = 1081.25
$ws.Range("K67").Value = 1452.5
$ws.Range("L67").Value = 1081.25
$ws.Range("M67").Value = -672.5
$ws.Range("N67").Value = -2641.25

$ws.Range("H80").Value = 1240.9615
$ws.Range("I80").Value = 1104.5
$ws.Range("J80").Value = 1400.1666
$ws.Range("K80").Value = 1104.5
$ws.Range("L80").Value = 1400.1666
$ws.Range("M80").Value = -106.5
$ws.Range("N80").Value = -3396.1666

$ws.Range("H83").Value = 1240.9615
$ws.Range("I83").Value = 1104.5
$ws.Range("J83").Value = 1400.1666
$ws.Range("K83").Value = 5522.5
$ws.Range("L83").Value = 7000.833000000001
$ws.Range("M83").Value = -530.5
$ws.Range("N83").Value = -16984.833

$ws.Range("H107").Value = 32327.412
$ws.Range("I107").Value = 57334.332
$ws.Range("J107").Value = 4194.625
$ws.Range("K107").Value = 57334.332
$ws.Range("L107").Value = 4194.625
$ws.Range("M107").Value = -55414.332
$ws.Range("N107").Value = -8034.625

$ws.Range("H134").Value = 2248.04
$ws.Range("I134").Value = 2032.1818
$ws.Range("J134").Value = 3831
$ws.Range("K134").Value = 6096.5454
$ws.Range("L134").Value = 11493
$ws.Range("M134").Value = -3561.5454
$ws.Range("N134").Value = -16563

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6251689.5
$ws.Range("I31").Value = 6668135.5
$ws.Range("K31").Value = 6668135.5
$ws.Range("M31").Value = -6667840.5

$ws.Range("H34").Value = 6251689.5
$ws.Range("I34").Value = 6668135.5
$ws.Range("K34").Value = 6668135.5
$ws.Range("M34").Value = -6667933.5

$ws.Range("H58").Value = 998.9259
$ws.Range("I58").Value = 1000.3043
$ws.Range("J58").Value = 991
$ws.Range("K58").Value = 1000.3043
$ws.Range("L58").Value = 991
$ws.Range("M58").Value = -797.3043
$ws.Range("N58").Value = -1397

$ws.Range("H105").Value = 1681.8572
$ws.Range("I105").Value = 1069.5555
$ws.Range("J105").Value = 2784
$ws.Range("K105").Value = 1069.5555
$ws.Range("L105").Value = 2784
$ws.Range("M105").Value = 677.4445000000001
$ws.Range("N105").Value = -6278

$ws.Range("H134").Value = 1390.4706
$ws.Range("I134").Value = 1149.9
$ws.Range("J134").Value = 3194.75
$ws.Range("K134").Value = 3449.7
$ws.Range("L134").Value = 9584.25
$ws.Range("M134").Value = -914.7000000000003
$ws.Range("N134").Value = -14654.25

$ws.Range("H136").Value = 998.9259
$ws.Range("I136").Value = 1000.3043
$ws.Range("J136").Value = 991
$ws.Range("K136").Value = 3000.9129
$ws.Range("L136").Value = 2973
$ws.Range("M136").Value = -450.9129000000003
$ws.Range("N136").Value = -8073

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 4713.75
$ws.Range("I112").Value = 4713.75
$ws.Range("K112").Value = 14141.25
$ws.Range("M112").Value = -13033.25

$ws.Range("H118").Value = 6333
$ws.Range("I118").Value = 6333
$ws.Range("K118").Value = 18999
$ws.Range("M118").Value = -17756

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3827.25
$ws.Range("I70").Value = 3991.52
$ws.Range("J70").Value = 3240.5715
$ws.Range("K70").Value = 3991.52
$ws.Range("L70").Value = 3240.5715
$ws.Range("M70").Value = -3721.52
$ws.Range("N70").Value = -3780.5715

$ws.Range("H73").Value = 3827.25
$ws.Range("I73").Value = 3991.52
$ws.Range("J73").Value = 3240.5715
$ws.Range("K73").Value = 3991.52
$ws.Range("L73").Value = 3240.5715
$ws.Range("M73").Value = -3055.52
$ws.Range("N73").Value = -5112.5715

$ws.Range("H102").Value = 26928.05
$ws.Range("I102").Value = 26928.05
$ws.Range("K102").Value = 26928.05
$ws.Range("M102").Value = -25306.05

$ws.Range("H132").Value = 1730.75
$ws.Range("I132").Value = 1697.0714
$ws.Range("J132").Value = 1966.5
$ws.Range("K132").Value = 5091.2142
$ws.Range("L132").Value = 5899.5
$ws.Range("M132").Value = -2561.2142
$ws.Range("N132").Value = -10959.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2803.1538
$ws.Range("I132").Value = 2940.6572
$ws.Range("K132").Value = 8821.971600000001
$ws.Range("M132").Value = -6291.971600000001

$ws.Range("H136").Value = 3086.6843
$ws.Range("I136").Value = 2758.9092
$ws.Range("J136").Value = 5250
$ws.Range("K136").Value = 8276.7276
$ws.Range("L136").Value = 15750
$ws.Range("M136").Value = -5726.7276
$ws.Range("N136").Value = -20850

$ws.Range("H137").Value = 40000
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4824466
$ws.Range("I132").Value = 5700725.5
$ws.Range("J132").Value = 5039.125
$ws.Range("K132").Value = 17102176.5
$ws.Range("L132").Value = 15117.375
$ws.Range("M132").Value = -17099646.5
$ws.Range("N132").Value = -20177.375

$ws.Range("H136").Value = 20302.5
$ws.Range("I136").Value = 22545.758
$ws.Range("J136").Value = 5497
$ws.Range("K136").Value = 67637.274
$ws.Range("L136").Value = 16491
$ws.Range("M136").Value = -65087.274
$ws.Range("N136").Value = -21591
